$wb = $excel.ActiveWorkbook

# ---- Sheet: ALC ----
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H11").Value = 53505.95
$ws.Range("I11").Value = 53505.95
$ws.Range("K11").Value = 53505.95
$ws.Range("M11").Value = -53365.95
$ws.Range("H31").Value = 1072.5
$ws.Range("I31").Value = 1072.5
$ws.Range("K31").Value = 3217.5
$ws.Range("M31").Value = -2987.5
$ws.Range("H48").Value = 3980
$ws.Range("I48").Value = 3980
$ws.Range("J48").Value = 0
$ws.Range("K48").Value = 11940
$ws.Range("L48").Value = 0
$ws.Range("M48").Value = -11648
$ws.Range("N48").ClearContents()
$ws.Range("H56").Value = 3980
$ws.Range("I56").Value = 3980
$ws.Range("J56").Value = 0
$ws.Range("K56").Value = 11940
$ws.Range("L56").Value = 0
$ws.Range("M56").Value = -11406
$ws.Range("N56").ClearContents()
$ws.Range("H64").Value = 3786.4285
$ws.Range("I64").Value = 3374.75
$ws.Range("K64").Value = 3374.75
$ws.Range("M64").Value = -3126.75
$ws.Range("H67").Value = 3786.4285
$ws.Range("I67").Value = 3374.75
$ws.Range("K67").Value = 3374.75
$ws.Range("M67").Value = -2516.75
$ws.Range("H98").Value = 2661.25
$ws.Range("I98").Value = 2770
$ws.Range("J98").Value = 1900
$ws.Range("K98").Value = 2770
$ws.Range("L98").Value = 1900
$ws.Range("M98").Value = -1272
$ws.Range("N98").Value = -4896
$ws.Range("H122").Value = 2661.25
$ws.Range("I122").Value = 2770
$ws.Range("J122").Value = 1900
$ws.Range("K122").Value = 8310
$ws.Range("L122").Value = 5700
$ws.Range("M122").Value = -5860
$ws.Range("N122").Value = -10600
$ws.Range("H138").Value = 2681
$ws.Range("I138").Value = 1650.9667
$ws.Range("J138").Value = 3473.3333
$ws.Range("K138").Value = 4952.9001
$ws.Range("L138").Value = 10419.9999
$ws.Range("M138").Value = 187.0999000000002
$ws.Range("N138").Value = -20699.9999

# ---- Sheet: ARM ----
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 1203.25
$ws.Range("I2").Value = 1159
$ws.Range("K2").Value = 1159
$ws.Range("M2").Value = -1046
$ws.Range("H21").Value = 66944.664
$ws.Range("I21").Value = 800
$ws.Range("J21").Value = 100017
$ws.Range("K21").Value = 800
$ws.Range("L21").Value = 100017
$ws.Range("M21").Value = -426
$ws.Range("N21").Value = -100765
$ws.Range("H45").Value = 1130.1428
$ws.Range("I45").Value = 985.1667
$ws.Range("J45").Value = 2000
$ws.Range("K45").Value = 985.1667
$ws.Range("L45").Value = 2000
$ws.Range("M45").Value = -608.1667
$ws.Range("N45").Value = -2754
$ws.Range("H109").Value = 16401.428
$ws.Range("J109").Value = 16401.428
$ws.Range("L109").Value = 16401.428
$ws.Range("N109").Value = -19175.428
$ws.Range("H110").Value = 1193.4762
$ws.Range("I110").Value = 812.86664
$ws.Range("J110").Value = 2145
$ws.Range("K110").Value = 812.86664
$ws.Range("L110").Value = 2145
$ws.Range("M110").Value = 1232.13336
$ws.Range("N110").Value = -6235
$ws.Range("H116").Value = 1203.25
$ws.Range("I116").Value = 1159
$ws.Range("K116").Value = 1159
$ws.Range("M116").Value = 1135
$ws.Range("H122").Value = 1510.1052
$ws.Range("I122").Value = 1739
$ws.Range("J122").Value = 1343.6364
$ws.Range("K122").Value = 5217
$ws.Range("L122").Value = 4030.9092
$ws.Range("M122").Value = -2767
$ws.Range("N122").Value = -8930.9092
$ws.Range("H125").Value = 47486
$ws.Range("J125").Value = 50928.75
$ws.Range("L125").Value = 50928.75
$ws.Range("N125").Value = -60768.75
$ws.Range("H134").Value = 59656
$ws.Range("J134").Value = 59656
$ws.Range("L134").Value = 59656
$ws.Range("N134").Value = -69796
$ws.Range("H135").Value = 100020470
$ws.Range("J135").Value = 100020470
$ws.Range("L135").Value = 100020470
$ws.Range("N135").Value = -100030610

# ---- Sheet: BSM ----
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 1203.25
$ws.Range("I3").Value = 1159
$ws.Range("K3").Value = 1159
$ws.Range("M3").Value = -1045
$ws.Range("H24").Value = 2003
$ws.Range("I24").Value = 2003
$ws.Range("K24").Value = 2003
$ws.Range("M24").Value = -1768
$ws.Range("H43").Value = 231886.67
$ws.Range("J43").Value = 231886.67
$ws.Range("L43").Value = 231886.67
$ws.Range("N43").Value = -232248.67
$ws.Range("H62").Value = 40000
$ws.Range("J62").Value = 40000
$ws.Range("L62").Value = 40000
$ws.Range("N62").Value = -41372
$ws.Range("H65").Value = 40000
$ws.Range("J65").Value = 40000
$ws.Range("L65").Value = 120000
$ws.Range("N65").Value = -126864
$ws.Range("H105").Value = 2313.077
$ws.Range("I105").Value = 2099.0908
$ws.Range("K105").Value = 2099.0908
$ws.Range("M105").Value = -352.0907999999999
$ws.Range("H107").Value = 2014.8334
$ws.Range("I107").Value = 2014.8334
$ws.Range("J107").Value = 0
$ws.Range("K107").Value = 2014.8334
$ws.Range("L107").Value = 0
$ws.Range("M107").Value = -94.83339999999998
$ws.Range("N107").ClearContents()

# ---- Sheet: CRP ----
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H2").Value = 964.3333
$ws.Range("I2").Value = 964.3333
$ws.Range("K2").Value = 964.3333
$ws.Range("M2").Value = -851.3333
$ws.Range("H5").Value = 344.17648
$ws.Range("I5").Value = 145.4
$ws.Range("J5").Value = 628.1429000000001
$ws.Range("K5").Value = 145.4
$ws.Range("L5").Value = 628.1429000000001
$ws.Range("M5").Value = -33.40000000000001
$ws.Range("N5").Value = -852.1429000000001
$ws.Range("H36").Value = 2497.25
$ws.Range("I36").Value = 1663
$ws.Range("K36").Value = 1663
$ws.Range("M36").Value = -1275
$ws.Range("H40").Value = 2497.25
$ws.Range("I40").Value = 1663
$ws.Range("K40").Value = 1663
$ws.Range("M40").Value = -1503
$ws.Range("H58").Value = 2936.7827
$ws.Range("I58").Value = 1370.7142
$ws.Range("J58").Value = 5372.8887
$ws.Range("K58").Value = 1370.7142
$ws.Range("L58").Value = 5372.8887
$ws.Range("M58").Value = -1167.7142
$ws.Range("N58").Value = -5778.8887
$ws.Range("H122").Value = 45455620
$ws.Range("I122").Value = 76923750
$ws.Range("K122").Value = 230771250
$ws.Range("M122").Value = -230768800
$ws.Range("H131").Value = 22326
$ws.Range("J131").Value = 22326
$ws.Range("L131").Value = 22326
$ws.Range("N131").Value = -32406
$ws.Range("H136").Value = 2936.7827
$ws.Range("I136").Value = 1370.7142
$ws.Range("J136").Value = 5372.8887
$ws.Range("K136").Value = 4112.142599999999
$ws.Range("L136").Value = 16118.6661
$ws.Range("M136").Value = -1562.142599999999
$ws.Range("N136").Value = -21218.6661

# ---- Sheet: GSM ----
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H9").Value = 1170.4286
$ws.Range("I9").Value = 697.5
$ws.Range("J9").Value = 4008
$ws.Range("K9").Value = 697.5
$ws.Range("L9").Value = 4008
$ws.Range("M9").Value = -527.5
$ws.Range("N9").Value = -4348
$ws.Range("H13").Value = 326.25
$ws.Range("J13").Value = 500
$ws.Range("L13").Value = 500
$ws.Range("N13").Value = -778
$ws.Range("H102").Value = 5106269.5
$ws.Range("I102").Value = 7146416
$ws.Range("K102").Value = 7146416
$ws.Range("M102").Value = -7144794
$ws.Range("H122").Value = 3815.923
$ws.Range("I122").Value = 3325.875
$ws.Range("J122").Value = 4600
$ws.Range("K122").Value = 9977.625
$ws.Range("L122").Value = 13800
$ws.Range("M122").Value = -7527.625
$ws.Range("N122").Value = -18700
$ws.Range("H123").Value = 19735.715
$ws.Range("J123").Value = 19735.715
$ws.Range("L123").Value = 19735.715
$ws.Range("N123").Value = -24635.715
$ws.Range("H134").Value = 24255.77
$ws.Range("J134").Value = 24255.77
$ws.Range("L134").Value = 72767.31
$ws.Range("N134").Value = -77837.31

# ---- Sheet: LTW ----
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H5").Value = 0
$ws.Range("J5").Value = 0
$ws.Range("L5").Value = 0
$ws.Range("N5").ClearContents()
$ws.Range("H7").Value = 4500
$ws.Range("J7").Value = 4500
$ws.Range("L7").Value = 4500
$ws.Range("N7").Value = -4724
$ws.Range("H126").Value = 4500
$ws.Range("J126").Value = 4500
$ws.Range("L126").Value = 13500
$ws.Range("N126").Value = -18440
$ws.Range("H134").Value = 48666.668
$ws.Range("J134").Value = 48666.668
$ws.Range("L134").Value = 48666.668
$ws.Range("N134").Value = -58806.668
$ws.Range("H138").Value = 0
$ws.Range("J138").Value = 0
$ws.Range("L138").Value = 0
$ws.Range("N138").ClearContents()

# ---- Sheet: WVR ----
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H46").Value = 83171.5
$ws.Range("J46").Value = 83171.5
$ws.Range("L46").Value = 83171.5
$ws.Range("N46").Value = -83633.5
$ws.Range("H94").Value = 15330
$ws.Range("J94").Value = 15330
$ws.Range("L94").Value = 15330
$ws.Range("N94").Value = -17132
$ws.Range("H122").Value = 2520.2
$ws.Range("I122").Value = 2077.2307
$ws.Range("J122").Value = 3342.8572
$ws.Range("K122").Value = 6231.6921
$ws.Range("L122").Value = 10028.5716
$ws.Range("M122").Value = -3781.6921
$ws.Range("N122").Value = -14928.5716
$ws.Range("H134").Value = 83171.5
$ws.Range("J134").Value = 83171.5
$ws.Range("L134").Value = 249514.5
$ws.Range("N134").Value = -254584.5
$ws.Range("H137").Value = 42357.5
$ws.Range("J137").Value = 42357.5
$ws.Range("L137").Value = 42357.5
$ws.Range("N137").Value = -52557.5

Write-Host "Done updating Mandragora Profits sheets"